$d = $word.ActiveDocument

# New bold, numbered (numId=1) heading paragraph: "Text input test"
$p1xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr>' +
    '<w:pStyle w:val="Normal"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:bidi w:val="0"/>' +
    '<w:jc w:val="start"/>' +
    '<w:rPr><w:b/><w:bCs/></w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>' +
      '<w:b/><w:bCs/>' +
      '<w:i w:val="false"/>' +
      '<w:color w:val="000000"/>' +
      '<w:sz w:val="22"/>' +
    '</w:rPr>' +
    '<w:t>Text input test</w:t>' +
  '</w:r>' +
'</w:p>'

# New non-bold, non-numbered (numId=0) indented description paragraph
$p2xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr>' +
    '<w:pStyle w:val="Normal"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
    '<w:bidi w:val="0"/>' +
    '<w:ind w:hanging="0" w:start="720"/>' +
    '<w:jc w:val="start"/>' +
    '<w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>' +
      '<w:b w:val="false"/><w:bCs w:val="false"/>' +
      '<w:i w:val="false"/>' +
      '<w:color w:val="000000"/>' +
      '<w:sz w:val="22"/>' +
    '</w:rPr>' +
    '<w:t>Verifying correct words by reading and inputting them simulating manual keyboard input.</w:t>' +
  '</w:r>' +
'</w:p>'

$r1 = $d.Range($d.Content.End, $d.Content.End)
$r1.InsertXML($p1xml) | Out-Null

$r2 = $d.Range($d.Content.End, $d.Content.End)
$r2.InsertXML($p2xml) | Out-Null

Write-Host "Paragraphs after edit:" $d.Paragraphs.Count
